$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: the three runs that make up
#   "El informe debe guardarse en el " + "directorio ." + "/data y debe
#   llamarse discosCostosos.txt."
# (with <w:proofErr gramStart/gramEnd> markers around the middle run)
# become a single run with the full sentence and no proofErr markers.
# ---------------------------------------------------------------------
$rng1 = $d.Content
$rng1.Find.ClearFormatting()
$rng1.Find.Replacement.ClearFormatting()
$found1 = $rng1.Find.Execute(
    "El informe debe guardarse en el directorio ./data y debe llamarse discosCostosos.txt.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "El informe debe guardarse en el directorio ./data y debe llamarse discosCostosos.txt.",
    2)
Write-Output "change1 found=$found1"

# ---------------------------------------------------------------------
# Change 2: "INTEGRANTES:  Jhon Fredy Bolaños Chávez" ->
#           "INTEGRANTES:  JONATHAN CAMILO BURBANO PAZOS "
# "Jhon" is wrapped in <w:proofErr spellStart/spellEnd>; those two runs
# plus the trailing " Fredy Bolaños Chávez" run collapse into a single
# new run and both proofErr markers disappear, while the lone " " run
# right after "INTEGRANTES: " is left alone (just re-isolated into its
# own run again).
# ---------------------------------------------------------------------

# Include the leading space in the match/replacement so the engine's
# run-merge logic actually crosses (and so drops) the spellStart marker
# that otherwise clings to the untouched start of the match.
$rng2 = $d.Content
$rng2.Find.ClearFormatting()
$rng2.Find.Replacement.ClearFormatting()
$found2 = $rng2.Find.Execute(
    " Jhon Fredy Bolaños Chávez",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    " JONATHAN CAMILO BURBANO PAZOS ",
    2)
Write-Output "change2 found=$found2"

# The line above also swallows the standalone " " run that sits between
# "INTEGRANTES: " and the name (it keeps the same formatting, so the
# merge is invisible in the rendered text). Re-isolate that leading
# space back into its own run so the paragraph's run layout matches the
# original structure (three runs: "INTEGRANTES: ", " ", the new name).
$rng3 = $d.Content
$rng3.Find.ClearFormatting()
$found3 = $rng3.Find.Execute("INTEGRANTES: ", $false)
Write-Output "anchor found=$found3"
$splitPos = $rng3.End
$spaceChar = $d.Range($splitPos, $splitPos + 1)
Write-Output "space char=[$($spaceChar.Text)]"
$wasBold = $spaceChar.Bold
$spaceChar.Bold = 1
$spaceChar.Bold = $wasBold
